$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = 2054412.45
$ws.Range("C7").Value = -54.5225672500711
$ws.Range("D7").Value = 1960
$ws.Range("E7").Value = 1960
$ws.Range("F7").Value = 1048.169617346939
$ws.Range("G7").Value = 8.171322183759443
